# Insert a new weekly data row after row 27 (pushing rows 28-57 down to 29-58)
# and populate it, then update the date in row 27 to reflect the new week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original row 27 values before we touch anything, since the
# new row 28 should carry the same data that row 27 used to have (minus
# the date update applied to row 27 itself).
$origDate27 = $ws.Range("D27").Value()
$origJ27 = $ws.Range("J27").Value()
$origK27 = $ws.Range("K27").Value()
$origL27 = $ws.Range("L27").Value()
$origM27 = $ws.Range("M27").Value()
$origN27 = $ws.Range("N27").Value()
$origO27 = $ws.Range("O27").Value()
$origP27 = $ws.Range("P27").Value()
$origQ27 = $ws.Range("Q27").Value()
$origR27 = $ws.Range("R27").Value()
$origA27 = $ws.Range("A27").Value()
$origB27 = $ws.Range("B27").Value()
$origC27 = $ws.Range("C27").Value()
$origE27 = $ws.Range("E27").Value()
$origF27 = $ws.Range("F27").Value()
$origG27 = $ws.Range("G27").Value()
$origH27 = $ws.Range("H27").Value()
$origI27 = $ws.Range("I27").Value()

# Insert a new blank row at position 28, shifting rows 28..57 down to 29..58.
$ws.Rows.Item(28).Insert()

# Row 27 now gets an updated (newer) date.
$ws.Range("D27").Value = 45049

# The newly inserted row 28 gets the data that used to live in row 27
# (same market/category/price info, with the original row 27 date).
$ws.Range("A28").Value = $origA27
$ws.Range("B28").Value = $origB27
$ws.Range("C28").Value = $origC27
$ws.Range("D28").Value = $origDate27
$ws.Range("E28").Value = $origE27
$ws.Range("F28").Value = $origF27
$ws.Range("G28").Value = $origG27
$ws.Range("H28").Value = $origH27
$ws.Range("I28").Value = $origI27
$ws.Range("J28").Value = $origJ27
$ws.Range("K28").Value = $origK27
$ws.Range("L28").Value = $origL27
$ws.Range("M28").Value = $origM27
$ws.Range("N28").Value = $origN27
$ws.Range("O28").Value = $origO27
$ws.Range("P28").Value = $origP27
$ws.Range("Q28").Value = $origQ27
$ws.Range("R28").Value = $origR27
